$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 978
$ws.Range("F6").Value = 5407
$ws.Range("F9").Value = 936
$ws.Range("F11").Value = 76
$ws.Range("F15").Value = 19
$ws.Range("F17").Value = 1804
$ws.Range("F22").Value = 324
$ws.Range("F23").Value = 532
$ws.Range("F28").Value = 2759
$ws.Range("F32").Value = 109
$ws.Range("F34").Value = 348
$ws.Range("F35").Value = 13
$ws.Range("F40").Value = 676

# Sheet "演出" (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 174
$ws.Range("F6").Value = 123

# Sheet "全部类型" (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 978
$ws.Range("F7").Value = 5407
$ws.Range("F11").Value = 174
$ws.Range("F12").Value = 936
$ws.Range("F15").Value = 123
$ws.Range("F16").Value = 76
$ws.Range("F20").Value = 19
$ws.Range("F23").Value = 1804
$ws.Range("F27").Value = 324
$ws.Range("F29").Value = 532
$ws.Range("F32").Value = 2759
$ws.Range("F36").Value = 109
$ws.Range("F38").Value = 348
$ws.Range("F39").Value = 13
$ws.Range("F43").Value = 676
